# Updated cryptos list on Sat Aug 26 12:51:54 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal-text value into a cell without Excel silently
# re-interpreting numeric-looking text (e.g. "217.35", "0.5270") as a
# Number and dropping formatting like trailing zeros. A leading apostrophe
# forces text entry; re-applying the original Style afterwards keeps the
# cell's formatting identical to how it started (no stray quote-prefix flag).
function Set-CellText($range, $text) {
    $cell = $ws.Range($range)
    $savedStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $savedStyle
}

Set-CellText 'D2' "26.069.17"
Set-CellText 'E2' "  -0.50%  "
Set-CellText 'D3' "1.650.51"
Set-CellText 'D5' "217.35"
Set-CellText 'E5' "  +0.07%  "
Set-CellText 'D6' "0.5270"
Set-CellText 'E6' "  +1.51%  "
Set-CellText 'E8' "  -1.67%  "
Set-CellText 'D9' "0.06307"
Set-CellText 'E9' "  +0.63%  "
Set-CellText 'D10' "20.31"
Set-CellText 'E10' "  -2.26%  "
Set-CellText 'D11' "0.07794"
Set-CellText 'E11' "  +0.02%  "
Set-CellText 'D12' "4.514"
Set-CellText 'E12' "  +0.97%  "
Set-CellText 'D13' "1.674.80"
Set-CellText 'E13' "  +0.74%  "
Set-CellText 'D14' "1.876.28"
Set-CellText 'E14' "  -0.57%  "
Set-CellText 'D15' "0.5474"
Set-CellText 'E15' "  +0.18%  "
Set-CellText 'D16' "0.0₅8179"
Set-CellText 'D17' "65.44"
Set-CellText 'D18' "26.052.61"
Set-CellText 'E18' "  -0.58%  "
Set-CellText 'E19' "  -0.16%  "
Set-CellText 'D20' "4.568"
Set-CellText 'E20' "  -0.85%  "
Set-CellText 'D21' "190.32"
Set-CellText 'E21' "  -0.80%  "
Set-CellText 'D22' "10.07"
Set-CellText 'E22' "  +0.19%  "
Set-CellText 'D23' "6.011"
Set-CellText 'E23' "  +0.18%  "
Set-CellText 'E24' "  -0.12%  "
Set-CellText 'D25' "143.48"
Set-CellText 'E25' "  +3.17%  "
Set-CellText 'D26' "0.1232"
Set-CellText 'E26' "  +0.78%  "
Set-CellText 'D27' "7.213"
Set-CellText 'E27' "  -0.96%  "
Set-CellText 'D28' "15.97"
Set-CellText 'E28' "  -0.96%  "
Set-CellText 'D29' "1.436"
Set-CellText 'E29' "  -0.19%  "
Set-CellText 'D30' "0.05800"
Set-CellText 'E30' "  -2.40%  "
Set-CellText 'E31' "  -0.24%  "
Set-CellText 'D32' "3.544"
Set-CellText 'E32' "  -0.11%  "
Set-CellText 'D33' "3.260"
Set-CellText 'E33' "  -0.19%  "
Set-CellText 'D34' "1.591"
Set-CellText 'E34' "  +0.45%  "
Set-CellText 'D35' "2.794"
Set-CellText 'E35' "  +0.91%  "
Set-CellText 'D36' "2.413"
Set-CellText 'E36' "  -0.30%  "
Set-CellText 'D37' "0.9415"
Set-CellText 'E37' "  -1.88%  "
Set-CellText 'D38' "0.5741"
Set-CellText 'E38' "  +0.91%  "
Set-CellText 'D39' "0.01602"
Set-CellText 'E39' "  +0.65%  "
Set-CellText 'B40' "Quant"
Set-CellText 'C40' "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-CellText 'D40' "105.16"
Set-CellText 'E40' "  +4.73%  "
Set-CellText 'B41' "TrustWalletToken"
Set-CellText 'C41' "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-CellText 'D41' "0.8492"
Set-CellText 'E41' "  +0.05%  "
Set-CellText 'E42' "  -0.11%  "
Set-CellText 'D43' "5.708"
Set-CellText 'E43' "  -4.68%  "
Set-CellText 'D44' "1.027.62"
Set-CellText 'E44' "  +2.66%  "
Set-CellText 'D45' "1.794.00"
Set-CellText 'E45' "  -0.41%  "
Set-CellText 'D46' "57.10"
Set-CellText 'E46' "  +1.13%  "
Set-CellText 'D47' "0.9999"
Set-CellText 'D48' "0.4330"
Set-CellText 'E48' "  -0.15%  "
Set-CellText 'B49' "Cronos"
Set-CellText 'C49' "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText 'D49' "0.05140"
Set-CellText 'E49' "  -0.39%  "
Set-CellText 'B50' "EnergySwap"
Set-CellText 'C50' "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText 'D50' "7.816"
Set-CellText 'E50' "  -2.86%  "
Set-CellText 'D51' "1.450"
Set-CellText 'E51' "  -0.26%  "
